# Apply automatic data refresh update: rows 10<->11 and rows 12<->13 had their
# record data swapped (row position stays, underlying observation values move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# ---- Row 10 (becomes the former row 11 data; loses the "Ringhack" comment) ----
$ws.Range("A10").Value = 130894760
$ws.Range("B10").Value = 79243
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 406786
$ws.Range("R10").Value = 7010890
Set-TextCell "Y10" "2026-01-18"
Set-TextCell "Z10" "14:31"
Set-TextCell "AA10" "2026-01-18"
Set-TextCell "AB10" "14:31"
$ws.Range("AC10").ClearContents()

# ---- Row 11 (becomes the former row 10 data; gains the "Ringhack" comment) ----
$ws.Range("A11").Value = 130894767
$ws.Range("B11").Value = 57884
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("Q11").Value = 407194
$ws.Range("R11").Value = 7011100
Set-TextCell "Y11" "2026-01-21"
Set-TextCell "Z11" "12:26"
Set-TextCell "AA11" "2026-01-21"
Set-TextCell "AB11" "12:26"
$ws.Range("AC11").Value = "Ringhack"

# ---- Row 12 (becomes the former row 13 data) ----
$ws.Range("A12").Value = 130894782
$ws.Range("Q12").Value = 407192
$ws.Range("R12").Value = 7011093
Set-TextCell "Y12" "2026-01-24"
Set-TextCell "Z12" "15:03"
Set-TextCell "AA12" "2026-01-24"
Set-TextCell "AB12" "15:03"

# ---- Row 13 (becomes the former row 12 data) ----
$ws.Range("A13").Value = 130894766
$ws.Range("Q13").Value = 407194
$ws.Range("R13").Value = 7011099
Set-TextCell "Y13" "2026-01-21"
Set-TextCell "Z13" "12:22"
Set-TextCell "AA13" "2026-01-21"
Set-TextCell "AB13" "12:22"
